$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the dataset (RM 232 and SC 92) ---
# Row 26 = "RM 232" ; deleting it shifts everything below up by one row, so the row
# that used to be 28 ("SC 92") is now at row 27.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Update individual "imputed" values inside the RM block (rows 2-23) ---
$ws.Range("D2").Value = -13.5
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = -6.4
$ws.Range("E5").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("D12").Value = -14.1
$ws.Range("D14").Value = ""
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("D22").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = -7

# --- Update values inside the SC block, now shifted up two rows (26-33) ---
$ws.Range("E27").Value = ""    # SC 101
$ws.Range("E29").Value = -6.8  # SC 119
$ws.Range("B30").Value = -19.7 # SC 120
$ws.Range("D31").Value = -13.7 # SC 132
$ws.Range("B32").Value = ""    # SC 193
$ws.Range("D33").Value = -14.1 # SC 232
